$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear the "Public ID" values in the first 4 populated data rows (rows 2-5),
#    keeping their existing cell formatting/style (auto-increment id will be
#    generated at runtime instead of being baked into the sample data).
$ws.Range("A2:A5").ClearContents()

# 2) Rows 6-11 used to hold "Theme 5".."Theme 10" / "Description 5".."Description 10"
#    sample rows. These become blank template rows identical in style to the
#    already-blank rows 12-16, so copy that formatting over and clear the values.
$ws.Range("A12:C12").Copy()
$ws.Range("A6:C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A6:C11").ClearContents()

# 3) Remove the now unused trailing blank template rows 17-22, shrinking the
#    sample sheet down to 16 rows total.
$ws.Range("A17:H22").EntireRow.Delete()
